$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 890.6667
$ws.Range("I6").Value = 963.3333
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 2889.9999
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -2777.9999
$ws.Range("N6").Value = -2024
$ws.Range("H8").Value = 316
$ws.Range("I8").Value = 316
$ws.Range("K8").Value = 948
$ws.Range("M8").Value = -809
$ws.Range("H112").Value = 1451.5385
$ws.Range("J112").Value = 1642.7273
$ws.Range("L112").Value = 4928.1819
$ws.Range("N112").Value = -7144.1819
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H127").Value = 1284.4615
$ws.Range("I127").Value = 400
$ws.Range("J127").Value = 1445.2727
$ws.Range("K127").Value = 1200
$ws.Range("L127").Value = 4335.8181
$ws.Range("M127").Value = 3760
$ws.Range("N127").Value = -14255.8181
$ws.Range("H138").Value = 2187
$ws.Range("I138").Value = 745.38464
$ws.Range("J138").Value = 4159.737
$ws.Range("K138").Value = 2236.15392
$ws.Range("L138").Value = 12479.211
$ws.Range("M138").Value = 2903.84608
$ws.Range("N138").Value = -22759.211
$ws.Range("H141").Value = 10344.55
$ws.Range("I141").Value = 3618.1875
$ws.Range("J141").Value = 37250
$ws.Range("K141").Value = 10854.5625
$ws.Range("L141").Value = 111750
$ws.Range("M141").Value = -5674.5625
$ws.Range("N141").Value = -122110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1446.5938
$ws.Range("I61").Value = 1073.3478
$ws.Range("J61").Value = 2400.4443
$ws.Range("K61").Value = 1073.3478
$ws.Range("L61").Value = 2400.4443
$ws.Range("M61").Value = -861.3478
$ws.Range("N61").Value = -2824.4443
$ws.Range("H74").Value = 8931897
$ws.Range("I74").Value = 17860972
$ws.Range("J74").Value = 2821.5715
$ws.Range("K74").Value = 17860972
$ws.Range("L74").Value = 2821.5715
$ws.Range("M74").Value = -17860098
$ws.Range("N74").Value = -4569.5715
$ws.Range("H77").Value = 8931897
$ws.Range("I77").Value = 17860972
$ws.Range("J77").Value = 2821.5715
$ws.Range("K77").Value = 89304860
$ws.Range("L77").Value = 14107.8575
$ws.Range("M77").Value = -89300492
$ws.Range("N77").Value = -22843.8575
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344
$ws.Range("H105").Value = 42000
$ws.Range("J105").Value = 42000
$ws.Range("L105").Value = 42000
$ws.Range("N105").Value = -48988
$ws.Range("H121").Value = 28084
$ws.Range("J121").Value = 28084
$ws.Range("L121").Value = 28084
$ws.Range("N121").Value = -31578
$ws.Range("H123").Value = 45321.25
$ws.Range("J123").Value = 45321.25
$ws.Range("L123").Value = 45321.25
$ws.Range("N123").Value = -55121.25
$ws.Range("H132").Value = 1745.6066
$ws.Range("I132").Value = 1547.5927
$ws.Range("J132").Value = 3273.1428
$ws.Range("K132").Value = 4642.7781
$ws.Range("L132").Value = 9819.428400000001
$ws.Range("M132").Value = -2112.7781
$ws.Range("N132").Value = -14879.4284
$ws.Range("H136").Value = 1446.5938
$ws.Range("I136").Value = 1073.3478
$ws.Range("J136").Value = 2400.4443
$ws.Range("K136").Value = 3220.0434
$ws.Range("L136").Value = 7201.3329
$ws.Range("M136").Value = -670.0434
$ws.Range("N136").Value = -12301.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 21881
$ws.Range("J100").Value = 21881
$ws.Range("L100").Value = 21881
$ws.Range("N100").Value = -24045
$ws.Range("H134").Value = 1995.32
$ws.Range("I134").Value = 1574.3684
$ws.Range("K134").Value = 4723.1052
$ws.Range("M134").Value = -2188.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3408185.8
$ws.Range("I31").Value = 8376100
$ws.Range("J31").Value = 1616.1428
$ws.Range("K31").Value = 8376100
$ws.Range("L31").Value = 1616.1428
$ws.Range("M31").Value = -8375805
$ws.Range("N31").Value = -2206.1428
$ws.Range("H34").Value = 3408185.8
$ws.Range("I34").Value = 8376100
$ws.Range("J34").Value = 1616.1428
$ws.Range("K34").Value = 8376100
$ws.Range("L34").Value = 1616.1428
$ws.Range("M34").Value = -8375898
$ws.Range("N34").Value = -2020.1428
$ws.Range("H36").Value = 16517.666
$ws.Range("I36").Value = 19700
$ws.Range("J36").Value = 14926.5
$ws.Range("K36").Value = 19700
$ws.Range("L36").Value = 14926.5
$ws.Range("M36").Value = -19312
$ws.Range("N36").Value = -15702.5
$ws.Range("H40").Value = 16517.666
$ws.Range("I40").Value = 19700
$ws.Range("J40").Value = 14926.5
$ws.Range("K40").Value = 19700
$ws.Range("L40").Value = 14926.5
$ws.Range("M40").Value = -19540
$ws.Range("N40").Value = -15246.5
$ws.Range("H94").Value = 1019.125
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 993.2857
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 993.2857
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -1895.2857
$ws.Range("H96").Value = 14545.454
$ws.Range("I96").Value = 10000
$ws.Range("J96").Value = 15000
$ws.Range("K96").Value = 10000
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -7254
$ws.Range("N96").Value = -20492
$ws.Range("H106").Value = 37300
$ws.Range("J106").Value = 37300
$ws.Range("L106").Value = 37300
$ws.Range("N106").Value = -39824
$ws.Range("H132").Value = 1263.8644
$ws.Range("I132").Value = 910.2708
$ws.Range("J132").Value = 2806.818
$ws.Range("K132").Value = 2730.8124
$ws.Range("L132").Value = 8420.454000000002
$ws.Range("M132").Value = -200.8123999999998
$ws.Range("N132").Value = -13480.454
$ws.Range("H134").Value = 2874.7334
$ws.Range("I134").Value = 3012.077
$ws.Range("J134").Value = 1982
$ws.Range("K134").Value = 9036.231
$ws.Range("L134").Value = 5946
$ws.Range("M134").Value = -6501.231
$ws.Range("N134").Value = -11016

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86177.13
$ws.Range("I2").Value = 123824.375
$ws.Range("J2").Value = 126.28571
$ws.Range("K2").Value = 742946.25
$ws.Range("L2").Value = 757.71426
$ws.Range("M2").Value = -742833.25
$ws.Range("N2").Value = -983.71426
$ws.Range("H121").Value = 1821.2916
$ws.Range("I121").Value = 326.5
$ws.Range("J121").Value = 2319.5557
$ws.Range("K121").Value = 979.5
$ws.Range("L121").Value = 6958.6671
$ws.Range("M121").Value = 330.5
$ws.Range("N121").Value = -9578.667099999999
$ws.Range("H131").Value = 9009910
$ws.Range("J131").Value = 9804876
$ws.Range("L131").Value = 29414628
$ws.Range("N131").Value = -29424708

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H132").Value = 2500.377
$ws.Range("I132").Value = 2374.7334
$ws.Range("J132").Value = 2853.75
$ws.Range("K132").Value = 7124.2002
$ws.Range("L132").Value = 8561.25
$ws.Range("M132").Value = -4594.2002
$ws.Range("N132").Value = -13621.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H136").Value = 8175.6313
$ws.Range("I136").Value = 26471.4
$ws.Range("J136").Value = 1641.4286
$ws.Range("K136").Value = 79414.20000000001
$ws.Range("L136").Value = 4924.2858
$ws.Range("M136").Value = -76864.20000000001
$ws.Range("N136").Value = -10024.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H103").Value = 34500
$ws.Range("J103").Value = 34500
$ws.Range("L103").Value = 34500
$ws.Range("N103").Value = -36844
$ws.Range("H123").Value = 40140.605
$ws.Range("J123").Value = 40140.605
$ws.Range("L123").Value = 40140.605
$ws.Range("N123").Value = -49940.605
$ws.Range("H132").Value = 2201.1428
$ws.Range("I132").Value = 1349.48
$ws.Range("J132").Value = 4330.3
$ws.Range("K132").Value = 4048.44
$ws.Range("L132").Value = 12990.9
$ws.Range("M132").Value = -1518.44
$ws.Range("N132").Value = -18050.9
$ws.Range("H136").Value = 1064.2941
$ws.Range("I136").Value = 686.64514
$ws.Range("J136").Value = 4966.6665
$ws.Range("K136").Value = 2059.93542
$ws.Range("L136").Value = 14899.9995
$ws.Range("M136").Value = 490.0645800000002
$ws.Range("N136").Value = -19999.9995
